# update with Bill Center ballot
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# New ballot row for Bill Center
$row = 19

$ws.Cells.Item($row, 1).Value = "Bill Center"     # A19 voter name

# Votes ("x") for: Barry Bonds(C), Roger Clemens(D), Roy Halladay(E),
# Todd Helton(F), Edgar Martinez(I), Fred McGriff(J), Mike Mussina(K),
# Mariano Rivera(O), Omar Vizquel(T), Larry Walker(V)
$votedCols = @("C","D","E","F","I","J","K","O","T","V")
foreach ($col in $votedCols) {
    $ws.Range("$col$row").Value = "x"
}

$ws.Range("AK$row").Value = 10            # n_votes
$ws.Range("AL$row").Value = "DM"          # source
$ws.Range("AM$row").Value = 43439         # date serial (2018-12-05)

# Reuse the same date style as the cell above (AM18) instead of creating
# a brand new number format
$ws.Range("AM18").Copy()
$ws.Range("AM$row").PasteSpecial(-4122)   # xlPasteFormats

# Update the active sheet view/selection to match the saved state
$ws.Activate()
$ws.Range("E15").Select()
